# Append new Bank Deposit rows 26-33 (2025-12-05 entries), matching the
# source layout of the existing sheet ("Bank Deposits": A-G, row 1 header).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns A (Collection Date), D (Deposit Date) and G (Entry Date) hold
# dd-mm-yyyy / yyyy-mm-dd text that Excel would otherwise auto-convert to a
# date serial number, so format them as Text first (matches the existing
# rows 2:25, which already store these as literal text strings).
$ws.Range("A26:A33").NumberFormat = "@"
$ws.Range("D26:D29").NumberFormat = "@"
$ws.Range("G26:G33").NumberFormat = "@"

# Row 26
$ws.Range("A26").Value = "05-12-2025"
$ws.Range("B26").Value = "010965012-Medha Sub Division Office Coll."
$ws.Range("C26").Value = "Cash"
$ws.Range("D26").Value = "2025-12-05"
$ws.Range("E26").Value = 28180
# F26 (Remark) left blank
$ws.Range("G26").Value = "2025-12-23"

# Row 27
$ws.Range("A27").Value = "05-12-2025"
$ws.Range("B27").Value = "020965017-Kai Lalsingrao Shinde Gr.Bid.S.S.Pat.Ltd Kudal Br. Kudal"
$ws.Range("C27").Value = "Cash"
$ws.Range("D27").Value = "2025-12-05"
$ws.Range("E27").Value = 59730
$ws.Range("F27").Value = "रद्द झालेल्या दोन रिसिटची रक्कम रु 2020.00 एवढी बँकेमार्फत नजरचुकीने भरली गेल्याने, रु 2020.00  एवढी रक्कम दिनांक 06.12.2025 च्या एकूण कलेक्शन मधून कमी भरणार असल्याचे सांगण्यात आले. ( आज रोजीचे एकूण कलेक्शन 57110.00 असे आहे )"
$ws.Range("G27").Value = "2025-12-23"

# Row 28
$ws.Range("A28").Value = "05-12-2025"
$ws.Range("B28").Value = "020965018-Kai Lalsingrao Shinde Gr.Big.Sheti Sah.Pat.Ltd. Br. Medha"
$ws.Range("C28").Value = "Cash"
$ws.Range("D28").Value = "2025-12-05"
$ws.Range("E28").Value = 32740
# F28 (Remark) left blank
$ws.Range("G28").Value = "2025-12-23"

# Row 29
$ws.Range("A29").Value = "05-12-2025"
$ws.Range("B29").Value = "020965021-KAI.LALSINGRAO BAPUSO SHINDE SAH.PAT.LTD.,KUDAL, BR.KARAHAR"
$ws.Range("C29").Value = "Cash"
$ws.Range("D29").Value = "2025-12-05"
$ws.Range("E29").Value = 4020
# F29 (Remark) left blank
$ws.Range("G29").Value = "2025-12-23"

# Row 30
$ws.Range("A30").Value = "05-12-2025"
$ws.Range("B30").Value = "020965017-Kai Lalsingrao Shinde Gr.Bid.S.S.Pat.Ltd Kudal Br. Kudal"
$ws.Range("C30").Value = "Cash"
# D30 (Deposit Date) left blank
$ws.Range("E30").Value = 0
$ws.Range("F30").Value = "रद्द झालेल्या दोन रिसिटची रक्कम रु 2020.00 एवढी बँकेमार्फत नजरचुकीने भरली गेल्याने, रु 2020.00  एवढी रक्कम दिनांक 06.12.2025 च्या एकूण कलेक्शन मधून कमी भरणार असल्याचे सांगण्यात आले. ( आज रोजीचे एकूण कलेक्शन 57110.00 असे आहे )"
$ws.Range("G30").Value = "2025-12-23"

# Row 31
$ws.Range("A31").Value = "05-12-2025"
$ws.Range("B31").Value = "020965017-Kai Lalsingrao Shinde Gr.Bid.S.S.Pat.Ltd Kudal Br. Kudal"
$ws.Range("C31").Value = "Cheque"
# D31 (Deposit Date) left blank
$ws.Range("E31").Value = 0
$ws.Range("F31").Value = "रद्द झालेल्या दोन रिसिटची रक्कम रु 2020.00 एवढी बँकेमार्फत नजरचुकीने भरली गेल्याने, रु 2020.00  एवढी रक्कम दिनांक 06.12.2025 च्या एकूण कलेक्शन मधून कमी भरणार असल्याचे सांगण्यात आले. ( आज रोजीचे एकूण कलेक्शन 57110.00 असे आहे )"
$ws.Range("G31").Value = "2025-12-23"

# Row 32
$ws.Range("A32").Value = "05-12-2025"
$ws.Range("B32").Value = "020965017-Kai Lalsingrao Shinde Gr.Bid.S.S.Pat.Ltd Kudal Br. Kudal"
$ws.Range("C32").Value = "NEFT"
# D32 (Deposit Date) left blank
$ws.Range("E32").Value = 0
$ws.Range("F32").Value = "रद्द झालेल्या दोन रिसिटची रक्कम रु 2020.00 एवढी बँकेमार्फत नजरचुकीने भरली गेल्याने, रु 2020.00  एवढी रक्कम दिनांक 06.12.2025 च्या एकूण कलेक्शन मधून कमी भरणार असल्याचे सांगण्यात आले. ( आज रोजीचे एकूण कलेक्शन 57110.00 असे आहे )"
$ws.Range("G32").Value = "2025-12-23"

# Row 33
$ws.Range("A33").Value = "05-12-2025"
$ws.Range("B33").Value = "020965017-Kai Lalsingrao Shinde Gr.Bid.S.S.Pat.Ltd Kudal Br. Kudal"
$ws.Range("C33").Value = "Total"
# D33 (Deposit Date) left blank
$ws.Range("E33").Value = 0
$ws.Range("F33").Value = "रद्द झालेल्या दोन रिसिटची रक्कम रु 2020.00 एवढी बँकेमार्फत नजरचुकीने भरली गेल्याने, रु 2020.00  एवढी रक्कम दिनांक 06.12.2025 च्या एकूण कलेक्शन मधून कमी भरणार असल्याचे सांगण्यात आले. ( आज रोजीचे एकूण कलेक्शन 57110.00 असे आहे )"
$ws.Range("G33").Value = "2025-12-23"

